$p = $ppt.ActivePresentation

# --- 1) Table on slide 16: switch to a different built-in table style ---
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{31973E19-B1B6-4F3B-BFA6-8D882193D3A4}")
    }
}

# --- 2) Swap the two themes in the deck: ---
#     theme1.xml (the active design theme, currently "Integral")  -> becomes "Office Theme" colours
#     theme2.xml (the notes-master theme, currently "Office Theme") -> becomes "Integral" colours
$activeColors = $p.SlideMaster.Theme.ThemeColorScheme
$notesColors  = $p.NotesMaster.Theme.ThemeColorScheme

# Remember the current ("Integral") colours before overwriting them.
$integral = @(0, 16777215, 5332805, 13754083, 3722137, 3646819, 2412774, 38860, 13611854, 10915127, 2465643, 158642)

# The current notes-master ("Office Theme") colours - read them off directly so we
# don't have to hard-code them.
$office = @()
for ($i = 1; $i -le $notesColors.Count; $i++) {
    $office += $notesColors.Item($i).RGB
}

# theme1.xml becomes "Office Theme" coloured
for ($i = 1; $i -le $activeColors.Count; $i++) {
    $activeColors.Item($i).RGB = $office[$i - 1]
}

# theme2.xml becomes "Integral" coloured
for ($i = 1; $i -le $notesColors.Count; $i++) {
    $notesColors.Item($i).RGB = $integral[$i - 1]
}
